$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# "Mission" mini-table (columns E:G) is shifted down by rows so a new
# attribute ("Niveau" / "int") can be inserted right above it, and the
# selected cell moves from M14 to K19.
#
# Style-template cells (stable cells elsewhere on the sheet, untouched
# by this edit) used as a source for PasteSpecial(Formats) so every
# destination cell ends up with the exact same cell format (s index)
# as in the target workbook:
#   s=2 -> A1   s=3 -> B3   s=4 -> B2   s=5 -> C2   s=7 -> E2   s=8 -> A8
# ------------------------------------------------------------------

function Set-Fmt($destRef, $templateRef) {
    $ws.Range($templateRef).Copy() | Out-Null
    $ws.Range($destRef).PasteSpecial(-4122) | Out-Null
}

# --- Row 15: E15 empty(s=7), F15 = "Niveau" (new string), G15 = "int" ---
Set-Fmt "E15" "E2"
$ws.Range("E15").Value = $null

Set-Fmt "F15" "B3"
$ws.Range("F15").Value = "Niveau"

Set-Fmt "G15" "B3"
$ws.Range("G15").Value = "int"

# --- Row 16: E16 empty(s=8), F16 = "DateHeure", G16 = "DateHeure" ---
Set-Fmt "E16" "A8"
$ws.Range("E16").Value = $null

Set-Fmt "F16" "B3"
$ws.Range("F16").Value = "DateHeure"

Set-Fmt "G16" "B3"
$ws.Range("G16").Value = "DateHeure"

# --- Row 17: remove E17:G17 (previously the "Mission" header row); row
#     shrinks from the 19.2 header height down to the regular 15.6 ---
$ws.Range("E17:G17").Clear() | Out-Null
$ws.Rows("17").RowHeight = 15.6

# --- Row 18: remove E18:G18 (previously "Id" / "int") ---
$ws.Range("E18:G18").Clear() | Out-Null

# --- Row 19: becomes the "Mission" header (Mission / Propriétés / Type);
#     row grows from 15.6 up to the 19.2 header height ---
$ws.Rows("19").RowHeight = 19.2

Set-Fmt "E19" "A1"
$ws.Range("E19").Value = "Mission"

Set-Fmt "F19" "A1"
$ws.Range("F19").Value = "Propriétés"

Set-Fmt "G19" "A1"
$ws.Range("G19").Value = "Type"

# --- Row 20: F20/G20 become "Id" / "int" (header-style cells s=4 / s=5) ---
Set-Fmt "F20" "B2"
$ws.Range("F20").Value = "Id"

Set-Fmt "G20" "C2"
$ws.Range("G20").Value = "int"

# --- Row 21: F21/G21 become "Description" / "varchar" ---
Set-Fmt "F21" "B3"
$ws.Range("F21").Value = "Description"

Set-Fmt "G21" "B3"
$ws.Range("G21").Value = "varchar"

# --- Row 22: E22 style changes to s=7, F22/G22 become "DateHeureAttribution" / "DateHeure" ---
Set-Fmt "E22" "E2"
$ws.Range("E22").Value = $null

Set-Fmt "F22" "B3"
$ws.Range("F22").Value = "DateHeureAttribution"

Set-Fmt "G22" "B3"
$ws.Range("G22").Value = "DateHeure"

# --- Row 23: add E23 (empty, s=7), F23/G23 = "DateHeureValidation" / "DateHeure" ---
Set-Fmt "E23" "E2"
$ws.Range("E23").Value = $null

Set-Fmt "F23" "B3"
$ws.Range("F23").Value = "DateHeureValidation"

Set-Fmt "G23" "B3"
$ws.Range("G23").Value = "DateHeure"

# --- Row 24: add E24 (empty, s=8), F24/G24 = "Terminée" / "bool" ---
Set-Fmt "E24" "A8"
$ws.Range("E24").Value = $null

Set-Fmt "F24" "B3"
$ws.Range("F24").Value = "Terminée"

Set-Fmt "G24" "B3"
$ws.Range("G24").Value = "bool"

# --- Update the active selection shown in the sheet view ---
$ws.Range("K19").Select() | Out-Null
